$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare formatting for the new header cells by copying the
# --- existing green header style (shared by A1:H1, J1, K1) onto the
# --- cells that need to become headers: I1, L1, M1.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("L1").PasteSpecial(-4122)
$ws.Range("M1").PasteSpecial(-4122)

# --- Move the "ID Gudang" / "Gudang" block two columns to the right
# --- (J/K -> L/M), and the row2/row3 data along with it, before we
# --- overwrite the header text so nothing gets clobbered.
$ws.Range("L2").Value = $ws.Range("J2").Value()
$ws.Range("M2").Value = $ws.Range("K2").Value()
$ws.Range("L3").Value = $ws.Range("J3").Value()
$ws.Range("M3").Value = $ws.Range("K3").Value()

$ws.Range("J2").Clear()
$ws.Range("K2").Clear()
$ws.Range("J3").Clear()
$ws.Range("K3").Clear()

# --- Now set the header row text. Two new headers are inserted at
# --- the front (No.Polisi, Tanggal SJ), "No Surat Jalan" becomes
# --- "No.Surat Jalan", and everything else shifts one column right;
# --- a new "Tahun" header is added at I1; J1/L1 both read
# --- "ID Gudang" (K1 stays empty) and M1 reads "Gudang".
$ws.Range("A1").Value = "No.Polisi"
$ws.Range("B1").Value = "Tanggal SJ"
$ws.Range("C1").Value = "No.Surat Jalan"
$ws.Range("D1").Value = "No Sales Order"
$ws.Range("E1").Value = "No Mesin"
$ws.Range("F1").Value = "No Rangka"
$ws.Range("G1").Value = "Tipe"
$ws.Range("H1").Value = "Warna"
$ws.Range("I1").Value = "Tahun"
$ws.Range("J1").Value = "ID Gudang"
$ws.Range("K1").Clear()
$ws.Range("L1").Value = "ID Gudang"
$ws.Range("M1").Value = "Gudang"

# --- Match the author's final selection (activeCell moved from J1 to L1).
$ws.Range("L1").Select() | Out-Null
